$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New HIGH/LOW/CLOSE/LTP/VOL/9:25 CLOSE figures for the 30-min refresh
# (rows 2..17 => B..G), matching the latest algo run that now also
# tracks ONGC and BHEL.
$data = @(
    @(1482.25, 1451.35, 3025.6, 1486.15, 81, 1471.75),
    @(1573.9, 1556.1, 1572.75, 1573.1, 17, 1570.75),
    @(51670, 51346.400000000001, 51500, 51490.75, 19, 51475.9),
    @(110.8, 108.79, 110.11, 111.1, 1028, 110.41),
    @(835.4, 819.7, 831.9, 838.4, 91, 833.95),
    @(700.8, 690.3, 697, 703.65, 112, 696.1),
    @(1234.8, 1219.9000000000001, 1224, 1229.9000000000001, 159, 1221.45),
    @(967.25, 947.35, 960.55, 967.65, 42, 962.35),
    @(25290, 25113, 25270, 25265.200000000001, 107, 25172.85),
    @(3071.6, 2991, 3037.2, 3057.5, 215, 3001.6),
    @(815.35, 808.75, 812.8, 818.9, 354, 809.75),
    @(1208.3499999999999, 1186.1500000000001, 1197.3499999999999, 1204.75, 29, 1200.5999999999999),
    @(1139.0999999999999, 1064.5, 1110, 1106.0999999999999, 324, 1073.95),
    @(153.69999999999999, 151.30000000000001, 153.25, 154.08000000000001, 929, 153),
    @(4549.75, 4492.75, 4520.95, 4540.55, 29, 4503),
    @(3548.85, 3506.05, 3540, 3557.4, 25, 3523.1)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        # columns B..G are 2..7
        $ws.Cells.Item($r, $c + 2).Value = $rowValues[$c]
    }
}

# Move the active selection to E4, matching where the analyst left off
$ws.Range("E4").Select()
